# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# The "Valor Mora" amounts for period 1802 (row 16) and period 1710 (row 20)
# were swapped: 1802 now shows 8533, and 1710 now shows 32000.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F16").Value = 8533
$ws.Range("F20").Value = 32000
